# Slides modules 8-9 and Source Code modules 4-6
#
# Target slide: "5. Looping and Arrays" -> slide 7 (the slide whose
# top-level group "object 3" sits at off x=0,y=1828800 / ext 12181840x5029200
# and which ends with the "Slide Number Placeholder 26" shape).
#
# EMU -> Point conversion used throughout (PowerPoint COM works in points,
# OOXML stores EMU; 1 pt = 12700 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Change 1: nudge the big top-level group 10160 EMU to the right ---
$grp = $s.Shapes.Item(2)
$grp.Left = 10160 / 12700

# --- Change 2: add a new "40" textbox after the slide-number placeholder ---
$left   = 8534400 / 12700
$top    = 5337175 / 12700
$width  = 581025 / 12700
$height = 368300 / 12700

$txt = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$txt.Name = "Text Box 27"
$txt.TextFrame.WordWrap = -1
$txt.TextFrame.AutoSize = 1
$txt.TextFrame.TextRange.Text = "40"
$txt.Height = $height
$txt.Fill.Visible = $false
